$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Apply updated TPM-derived values (NATMI LR-pair stats) to the data rows (2-21).
# Each assignment below mirrors one changed cell from the source diff: numeric
# statistic recalculations plus the "Target cluster" text relabeling that
# resulted from the updated TPM input (shared-string swaps in column D).

$ws.Cells.Item(2, 7).Value = 4.573795666666666
$ws.Cells.Item(2, 8).Value = 13.721387
$ws.Cells.Item(2, 9).Value = 0.6529099782872276
$ws.Cells.Item(2, 10).Value = 0.6529099782872277
$ws.Cells.Item(2, 13).Value = 2.330737666666666
$ws.Cells.Item(2, 14).Value = 6.992213
$ws.Cells.Item(2, 15).Value = 0.08387223179237234
$ws.Cells.Item(2, 16).Value = 0.08387223179237237
$ws.Cells.Item(2, 17).Value = 10.66031783993678
$ws.Cells.Item(2, 18).Value = 95.942860559431
$ws.Cells.Item(2, 19).Value = 0.05476101703845915
$ws.Cells.Item(2, 20).Value = 0.05476101703845918
$ws.Cells.Item(3, 7).Value = 4.573795666666666
$ws.Cells.Item(3, 8).Value = 13.721387
$ws.Cells.Item(3, 9).Value = 0.6529099782872276
$ws.Cells.Item(3, 10).Value = 0.6529099782872277
$ws.Cells.Item(3, 15).Value = 0.8164642509041029
$ws.Cells.Item(3, 16).Value = 0.816464250904103
$ws.Cells.Item(3, 17).Value = 103.774136368876
$ws.Cells.Item(3, 18).Value = 933.9672273198839
$ws.Cells.Item(3, 19).Value = 0.5330776563300954
$ws.Cells.Item(3, 20).Value = 0.5330776563300955
$ws.Cells.Item(4, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(4, 7).Value = 4.573795666666666
$ws.Cells.Item(4, 8).Value = 13.721387
$ws.Cells.Item(4, 9).Value = 0.6529099782872276
$ws.Cells.Item(4, 10).Value = 0.6529099782872277
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.04840566666666667
$ws.Cells.Item(4, 14).Value = 0.145217
$ws.Cells.Item(4, 15).Value = 0.001741891141501687
$ws.Cells.Item(4, 16).Value = 0.001741891141501688
$ws.Cells.Item(4, 17).Value = 0.2213976284421111
$ws.Cells.Item(4, 18).Value = 1.992578655979
$ws.Cells.Item(4, 19).Value = 0.001137298107376581
$ws.Cells.Item(4, 20).Value = 0.001137298107376581
$ws.Cells.Item(5, 4).Value = "MuSCs"
$ws.Cells.Item(5, 7).Value = 4.573795666666666
$ws.Cells.Item(5, 8).Value = 13.721387
$ws.Cells.Item(5, 9).Value = 0.6529099782872276
$ws.Cells.Item(5, 10).Value = 0.6529099782872277
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 2.721158333333333
$ws.Cells.Item(5, 14).Value = 8.163475
$ws.Cells.Item(5, 15).Value = 0.09792162616202293
$ws.Cells.Item(5, 16).Value = 0.09792162616202295
$ws.Cells.Item(5, 17).Value = 12.44602219331389
$ws.Cells.Item(5, 18).Value = 112.014199739825
$ws.Cells.Item(5, 19).Value = 0.06393400681129642
$ws.Cells.Item(5, 20).Value = 0.06393400681129643
$ws.Cells.Item(6, 9).Value = 0.1985019229157801
$ws.Cells.Item(6, 10).Value = 0.1985019229157801
$ws.Cells.Item(6, 13).Value = 2.330737666666666
$ws.Cells.Item(6, 14).Value = 6.992213
$ws.Cells.Item(6, 15).Value = 0.08387223179237234
$ws.Cells.Item(6, 16).Value = 0.08387223179237237
$ws.Cells.Item(6, 17).Value = 3.241018916071666
$ws.Cells.Item(6, 18).Value = 29.169170244645
$ws.Cells.Item(6, 19).Value = 0.01664879929002394
$ws.Cells.Item(6, 20).Value = 0.01664879929002394
$ws.Cells.Item(7, 9).Value = 0.1985019229157801
$ws.Cells.Item(7, 10).Value = 0.1985019229157801
$ws.Cells.Item(7, 15).Value = 0.8164642509041029
$ws.Cells.Item(7, 16).Value = 0.816464250904103
$ws.Cells.Item(7, 19).Value = 0.1620697237964564
$ws.Cells.Item(7, 20).Value = 0.1620697237964564
$ws.Cells.Item(8, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(8, 9).Value = 0.1985019229157801
$ws.Cells.Item(8, 10).Value = 0.1985019229157801
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.04840566666666667
$ws.Cells.Item(8, 14).Value = 0.145217
$ws.Cells.Item(8, 15).Value = 0.001741891141501687
$ws.Cells.Item(8, 16).Value = 0.001741891141501688
$ws.Cells.Item(8, 17).Value = 0.06731074181166667
$ws.Cells.Item(8, 18).Value = 0.605796676305
$ws.Cells.Item(8, 19).Value = 0.0003457687410980482
$ws.Cells.Item(8, 20).Value = 0.0003457687410980483
$ws.Cells.Item(9, 4).Value = "MuSCs"
$ws.Cells.Item(9, 9).Value = 0.1985019229157801
$ws.Cells.Item(9, 10).Value = 0.1985019229157801
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 2.721158333333333
$ws.Cells.Item(9, 14).Value = 8.163475
$ws.Cells.Item(9, 15).Value = 0.09792162616202293
$ws.Cells.Item(9, 16).Value = 0.09792162616202295
$ws.Cells.Item(9, 17).Value = 3.783920326208333
$ws.Cells.Item(9, 18).Value = 34.055282935875
$ws.Cells.Item(9, 19).Value = 0.01943763108820172
$ws.Cells.Item(9, 20).Value = 0.01943763108820172
$ws.Cells.Item(10, 5).Value = 1
$ws.Cells.Item(10, 6).Value = 0.3333333333333333
$ws.Cells.Item(10, 7).Value = 0.3547236666666667
$ws.Cells.Item(10, 8).Value = 1.064171
$ws.Cells.Item(10, 9).Value = 0.05063685358512936
$ws.Cells.Item(10, 10).Value = 0.05063685358512936
$ws.Cells.Item(10, 13).Value = 2.330737666666666
$ws.Cells.Item(10, 14).Value = 6.992213
$ws.Cells.Item(10, 15).Value = 0.08387223179237234
$ws.Cells.Item(10, 16).Value = 0.08387223179237237
$ws.Cells.Item(10, 17).Value = 0.826767811158111
$ws.Cells.Item(10, 18).Value = 7.440910300422999
$ws.Cells.Item(10, 19).Value = 0.004247025921128389
$ws.Cells.Item(10, 20).Value = 0.004247025921128391
$ws.Cells.Item(11, 5).Value = 1
$ws.Cells.Item(11, 6).Value = 0.3333333333333333
$ws.Cells.Item(11, 7).Value = 0.3547236666666667
$ws.Cells.Item(11, 8).Value = 1.064171
$ws.Cells.Item(11, 9).Value = 0.05063685358512936
$ws.Cells.Item(11, 10).Value = 0.05063685358512936
$ws.Cells.Item(11, 15).Value = 0.8164642509041029
$ws.Cells.Item(11, 16).Value = 0.816464250904103
$ws.Cells.Item(11, 17).Value = 8.048269936107999
$ws.Cells.Item(11, 18).Value = 72.434429424972
$ws.Cells.Item(11, 19).Value = 0.04134318073052338
$ws.Cells.Item(11, 20).Value = 0.04134318073052339
$ws.Cells.Item(12, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(12, 5).Value = 1
$ws.Cells.Item(12, 6).Value = 0.3333333333333333
$ws.Cells.Item(12, 7).Value = 0.3547236666666667
$ws.Cells.Item(12, 8).Value = 1.064171
$ws.Cells.Item(12, 9).Value = 0.05063685358512936
$ws.Cells.Item(12, 10).Value = 0.05063685358512936
$ws.Cells.Item(12, 11).Value = 1
$ws.Cells.Item(12, 12).Value = 0.3333333333333333
$ws.Cells.Item(12, 13).Value = 0.04840566666666667
$ws.Cells.Item(12, 14).Value = 0.145217
$ws.Cells.Item(12, 15).Value = 0.001741891141501687
$ws.Cells.Item(12, 16).Value = 0.001741891141501688
$ws.Cells.Item(12, 17).Value = 0.01717063556744445
$ws.Cells.Item(12, 18).Value = 0.154535720107
$ws.Cells.Item(12, 19).Value = 0.00008820388669345478
$ws.Cells.Item(12, 20).Value = 0.00008820388669345479
$ws.Cells.Item(13, 4).Value = "MuSCs"
$ws.Cells.Item(13, 5).Value = 1
$ws.Cells.Item(13, 6).Value = 0.3333333333333333
$ws.Cells.Item(13, 7).Value = 0.3547236666666667
$ws.Cells.Item(13, 8).Value = 1.064171
$ws.Cells.Item(13, 9).Value = 0.05063685358512936
$ws.Cells.Item(13, 10).Value = 0.05063685358512936
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 2.721158333333333
$ws.Cells.Item(13, 14).Value = 8.163475
$ws.Cells.Item(13, 15).Value = 0.09792162616202293
$ws.Cells.Item(13, 16).Value = 0.09792162616202295
$ws.Cells.Item(13, 17).Value = 0.9652592615805555
$ws.Cells.Item(13, 18).Value = 8.687333354225
$ws.Cells.Item(13, 19).Value = 0.004958443046784128
$ws.Cells.Item(13, 20).Value = 0.004958443046784129
$ws.Cells.Item(14, 7).Value = 0.479723
$ws.Cells.Item(14, 8).Value = 1.439169
$ws.Cells.Item(14, 9).Value = 0.06848052609708123
$ws.Cells.Item(14, 10).Value = 0.06848052609708123
$ws.Cells.Item(14, 13).Value = 2.330737666666666
$ws.Cells.Item(14, 14).Value = 6.992213
$ws.Cells.Item(14, 15).Value = 0.08387223179237234
$ws.Cells.Item(14, 16).Value = 0.08387223179237237
$ws.Cells.Item(14, 17).Value = 1.118108465666333
$ws.Cells.Item(14, 18).Value = 10.062976190997
$ws.Cells.Item(14, 19).Value = 0.005743614558078
$ws.Cells.Item(14, 20).Value = 0.005743614558078002
$ws.Cells.Item(15, 7).Value = 0.479723
$ws.Cells.Item(15, 8).Value = 1.439169
$ws.Cells.Item(15, 9).Value = 0.06848052609708123
$ws.Cells.Item(15, 10).Value = 0.06848052609708123
$ws.Cells.Item(15, 15).Value = 0.8164642509041029
$ws.Cells.Item(15, 16).Value = 0.816464250904103
$ws.Cells.Item(15, 17).Value = 10.884360310212
$ws.Cells.Item(15, 18).Value = 97.95924279190798
$ws.Cells.Item(15, 19).Value = 0.0559119014413723
$ws.Cells.Item(15, 20).Value = 0.05591190144137231
$ws.Cells.Item(16, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(16, 7).Value = 0.479723
$ws.Cells.Item(16, 8).Value = 1.439169
$ws.Cells.Item(16, 9).Value = 0.06848052609708123
$ws.Cells.Item(16, 10).Value = 0.06848052609708123
$ws.Cells.Item(16, 11).Value = 1
$ws.Cells.Item(16, 12).Value = 0.3333333333333333
$ws.Cells.Item(16, 13).Value = 0.04840566666666667
$ws.Cells.Item(16, 14).Value = 0.145217
$ws.Cells.Item(16, 15).Value = 0.001741891141501687
$ws.Cells.Item(16, 16).Value = 0.001741891141501688
$ws.Cells.Item(16, 17).Value = 0.02322131163033334
$ws.Cells.Item(16, 18).Value = 0.208991804673
$ws.Cells.Item(16, 19).Value = 0.0001192856217738809
$ws.Cells.Item(16, 20).Value = 0.0001192856217738809
$ws.Cells.Item(17, 4).Value = "MuSCs"
$ws.Cells.Item(17, 7).Value = 0.479723
$ws.Cells.Item(17, 8).Value = 1.439169
$ws.Cells.Item(17, 9).Value = 0.06848052609708123
$ws.Cells.Item(17, 10).Value = 0.06848052609708123
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 2.721158333333333
$ws.Cells.Item(17, 14).Value = 8.163475
$ws.Cells.Item(17, 15).Value = 0.09792162616202293
$ws.Cells.Item(17, 16).Value = 0.09792162616202295
$ws.Cells.Item(17, 17).Value = 1.305402239141667
$ws.Cells.Item(17, 18).Value = 11.748620152275
$ws.Cells.Item(17, 19).Value = 0.006705724475857044
$ws.Cells.Item(17, 20).Value = 0.006705724475857044
$ws.Cells.Item(18, 5).Value = 2
$ws.Cells.Item(18, 6).Value = 0.6666666666666666
$ws.Cells.Item(18, 7).Value = 0.2064496666666667
$ws.Cells.Item(18, 8).Value = 0.619349
$ws.Cells.Item(18, 9).Value = 0.02947071911478163
$ws.Cells.Item(18, 10).Value = 0.02947071911478163
$ws.Cells.Item(18, 13).Value = 2.330737666666666
$ws.Cells.Item(18, 14).Value = 6.992213
$ws.Cells.Item(18, 15).Value = 0.08387223179237234
$ws.Cells.Item(18, 16).Value = 0.08387223179237237
$ws.Cells.Item(18, 17).Value = 0.4811800143707777
$ws.Cells.Item(18, 18).Value = 4.330620129337
$ws.Cells.Item(18, 19).Value = 0.002471774984682863
$ws.Cells.Item(18, 20).Value = 0.002471774984682864
$ws.Cells.Item(19, 5).Value = 2
$ws.Cells.Item(19, 6).Value = 0.6666666666666666
$ws.Cells.Item(19, 7).Value = 0.2064496666666667
$ws.Cells.Item(19, 8).Value = 0.619349
$ws.Cells.Item(19, 9).Value = 0.02947071911478163
$ws.Cells.Item(19, 10).Value = 0.02947071911478163
$ws.Cells.Item(19, 15).Value = 0.8164642509041029
$ws.Cells.Item(19, 16).Value = 0.816464250904103
$ws.Cells.Item(19, 17).Value = 4.684104280852
$ws.Cells.Item(19, 18).Value = 42.156938527668
$ws.Cells.Item(19, 19).Value = 0.02406178860565541
$ws.Cells.Item(19, 20).Value = 0.02406178860565542
$ws.Cells.Item(20, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(20, 5).Value = 2
$ws.Cells.Item(20, 6).Value = 0.6666666666666666
$ws.Cells.Item(20, 7).Value = 0.2064496666666667
$ws.Cells.Item(20, 8).Value = 0.619349
$ws.Cells.Item(20, 9).Value = 0.02947071911478163
$ws.Cells.Item(20, 10).Value = 0.02947071911478163
$ws.Cells.Item(20, 11).Value = 1
$ws.Cells.Item(20, 12).Value = 0.3333333333333333
$ws.Cells.Item(20, 13).Value = 0.04840566666666667
$ws.Cells.Item(20, 14).Value = 0.145217
$ws.Cells.Item(20, 15).Value = 0.001741891141501687
$ws.Cells.Item(20, 16).Value = 0.001741891141501688
$ws.Cells.Item(20, 17).Value = 0.009993333748111113
$ws.Cells.Item(20, 18).Value = 0.08994000373300001
$ws.Cells.Item(20, 19).Value = 0.00005133478455972257
$ws.Cells.Item(20, 20).Value = 0.00005133478455972258
$ws.Cells.Item(21, 4).Value = "MuSCs"
$ws.Cells.Item(21, 5).Value = 2
$ws.Cells.Item(21, 6).Value = 0.6666666666666666
$ws.Cells.Item(21, 7).Value = 0.2064496666666667
$ws.Cells.Item(21, 8).Value = 0.619349
$ws.Cells.Item(21, 9).Value = 0.02947071911478163
$ws.Cells.Item(21, 10).Value = 0.02947071911478163
$ws.Cells.Item(21, 11).Value = 3
$ws.Cells.Item(21, 12).Value = 1
$ws.Cells.Item(21, 13).Value = 2.721158333333333
$ws.Cells.Item(21, 14).Value = 8.163475
$ws.Cells.Item(21, 15).Value = 0.09792162616202293
$ws.Cells.Item(21, 16).Value = 0.09792162616202295
$ws.Cells.Item(21, 17).Value = 0.5617822308638889
$ws.Cells.Item(21, 18).Value = 5.056040077775
$ws.Cells.Item(21, 19).Value = 0.00288582073988363
$ws.Cells.Item(21, 20).Value = 0.002885820739883631
